# Update the "人民法院审理刑事案件罪犯情况" worksheet:
#  - drop the oldest decade (2000年-2009年) from the top of the table
#  - append a new row for 2021年 at the bottom

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 10 rows for years 2000年..2009年 (current rows 2-11).
# Deleting the whole block shifts 2010年..2020年 up to rows 2-12.
$ws.Range("A2:E11").EntireRow.Delete()

# Append the new 2021年 row right after the shifted data (now row 12),
# i.e. the new data goes into row 13.
$newRow = 13
$ws.Cells.Item($newRow, 1).Value = "2021年"
$ws.Cells.Item($newRow, 2).Value = 248949
$ws.Cells.Item($newRow, 3).Value = 34616
$ws.Cells.Item($newRow, 4).Value = 1714942
$ws.Cells.Item($newRow, 5).Value = 283565

# Match the styling used by the other year-label cells in column A
# (bold, centered, bordered) by copying the format from the row above.
$ws.Cells.Item($newRow - 1, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
